$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The alcohol measurement table had a duplicate/obsolete data column (M).
# Deleting it shifts the following column (old N) left to become the new M,
# shrinking the sheet's used range from A1:N119 to A1:M119.
$ws.Columns.Item(13).Delete()

# Reflect the new active cell position on the sheet after the column removal.
$ws.Range("M1").Select()
